# Add SoundManager functions - Mute function, Changing volume.
# The per-sound "Default Volume" (volume) column is no longer needed (muting /
# volume is now handled centrally by the SoundManager), so it is removed from
# both the SOUND_BUNDLE and SOUND_RESOURCE tables.

$wb = $excel.ActiveWorkbook

# --- SOUND_BUNDLE sheet: remove the "volume" column (column E) ---
$wsBundle = $wb.Worksheets.Item("SOUND_BUNDLE")
$wsBundle.Columns("E").Delete()

# --- SOUND_RESOURCE sheet: remove the "volume" column (column F) ---
$wsResource = $wb.Worksheets.Item("SOUND_RESOURCE")
$wsResource.Columns("F").Delete()

# SOUND_RESOURCE becomes the active/selected sheet.
$wsResource.Activate()
